# Update the "Setor_UO_Órgão" lookup table: append 4 new rows of government
# sector / agency reference data, extend the table & filter range to match,
# apply a black font color across the data, and widen column B.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Extend the Excel Table ("Tabela2") down to the new last row ------
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:E120"))

# --- 2. Write the four new rows of data -----------------------------------
$newRows = @(
    @(4731, "FUNDO DE DESENVOLVIMENTO DO MINISTÉRIO PÚBLICO", "MINISTÉRIO PÚBLICO", 1090, "PROCURADORIA GERAL DE JUSTIÇA"),
    @(4751, "FUNDO ESPECIAL DA ADVOCACIA GERAL DO ESTADO ", "ADVOCACIA GERAL", 1080, "ADVOCACIA GERAL DO ESTADO"),
    @(4741, "FUNDO ESPECIAL DE GARANTIA DE ACESSO À JUSTIÇA ", "DEFENSORIA PÚBLICA", 1440, "DEFENSORIA PUBLICA DO ESTADO DE MINAS GERAIS"),
    @(2471, "AGÊNCIA REGULADORA DE TRANSPORTES DO ESTADO DE MINAS GERAIS ", "INFRAESTRUTURA, MOBILIDADE E PARCERIAS", 1300, "SECRETARIA DE ESTADO DE INFRAESTRUTURA, MOBILIDADE E PARCERIAS")
)

$r = 117
foreach ($row in $newRows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $r = $r + 1
}

# --- 3. Update the hidden _FilterDatabase defined name to the new range ---
$fd = $wb.Names.Item(1)
$fd.RefersTo = "=Setor_UO_Orgao!`$A`$1:`$E`$120"

# --- 4. Apply a black font color across the whole table range -------------
$ws.Range("A1:E120").Font.Color = 0

# --- 5. Widen column B to fit the longer descriptions ----------------------
$ws.Columns("B:B").ColumnWidth = 62.14

# --- 6. Misc view / print setup housekeeping --------------------------------
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

$win = $excel.ActiveWindow
$win.ScrollRow = 108
$win.ScrollColumn = 1
$ws.Range("E107").Select()
